$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All tasks are now complete: set the status column (C) to "ολοκληρώθεκε" / completed
# for every data row (rows 2 through 10).
$status = "ολοκληρώθηκε"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = $status
}

# Move/update the active selection to C11 (was A11)
$ws.Range("C11").Select()
